$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Training Dashboard sheet: for data rows 3..35, the "LAST UPDATE" date (col I)
# moved forward from 08-Sep-2025 to 16-Sep-2025 (8 days later), so the
# "PERIOD TO EXPIRE" (col H) shrank by 8 for every row.
# ---------------------------------------------------------------------------
$wsTrain = $wb.Worksheets.Item("Training Dashboard")

# "LAST UPDATE" column (I) as a block: force text so "16-Sep-2025" is stored
# as a literal string instead of being auto-parsed into a date serial
# number, then restore the General number format so the cells' look is
# unchanged.
$dateRange = $wsTrain.Range("I3:I35")
$dateRange.NumberFormat = "@"
$dateRange.Value = "16-Sep-2025"
$dateRange.NumberFormat = "General"

for ($r = 3; $r -le 35; $r++) {
    $periodCell = $wsTrain.Cells.Item($r, 8)
    $oldPeriod = $periodCell.Value()
    $periodCell.Value = $oldPeriod - 8
}

# ---------------------------------------------------------------------------
# Exam Dashboard sheet: the date-validity checks were added/standardized so
# every comment now reads "date is valid"; the comments column is also
# narrower now that the text is uniform.
# ---------------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

$wsExam.Range("E3").Value = "date is valid"
$wsExam.Range("E4").Value = "date is valid"
$wsExam.Range("E5").Value = "date is valid"
$wsExam.Range("E6").Value = "date is valid"

## 15 characters of displayed column width; COM's ColumnWidth excludes the
## ~0.83-character cell-padding that ends up in the saved <col width=.../>,
## so back that padding out to land on a clean "15" in the saved file.
$wsExam.Columns.Item(5).ColumnWidth = 14.17

# ---------------------------------------------------------------------------
# Styling: the bold header font loses its fixed 14pt size and instead gets a
# white color; the section-header fill (row 2 on both sheets) now reuses that
# same bold/white font instead of its own plain-bold font.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsTrain, $wsExam)) {
    $titleRange = $ws.Range("A1")
    $titleRange.Font.Size = 11
    $titleRange.Font.Bold = $true
    $titleRange.Font.Color = 16777215

    $headerRange = $ws.Rows.Item(2)
    $headerRange.Font.Bold = $true
    $headerRange.Font.Color = 16777215
}

Write-Host "edit complete"
